# Weather forecast display, gameplay tweaks, primary level 1 plotting
# - Adds a new "weatherForecast" / "Weather Forecast" row (new row 12),
#   shifting the existing weather* rows down by one.
# - Adds a new "weatherLightRain" / "weatherLightRainDesc" pair of rows
#   right after "weatherClearDesc" (new rows 23-24).
# - Updates the active selection to reflect where editing left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the "Weather Forecast" heading row before weatherSunny (row 12) ---
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value = "weatherForecast"
$ws.Range("B12").Value = "Weather Forecast"

# --- Insert the new "Light Rain" weather entry after weatherClear/weatherClearDesc ---
# (weatherClearDesc is now at row 22 after the insertion above)
$ws.Rows.Item(23).Insert()
$ws.Rows.Item(24).Insert()

# Set key column first, then value column, to match shared-string insertion order
$ws.Range("A23").Value = "weatherLightRain"
$ws.Range("A24").Value = "weatherLightRainDesc"
$ws.Range("B23").Value = "Light Rain"
$ws.Range("B24").Value = "Drip here and there."

# --- Reflect the final cursor/selection position ---
$ws.Range("B24").Select()
